$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Table-of-Contents entries: merge the "Assignment NN" run with the
#    trailing whitespace run(s) that follow it into a single run.
#    (Find/Execute over the paragraph's own Range naturally folds the
#    multiple adjacent same-format runs into one run when replaced.)
# ---------------------------------------------------------------------------
function Merge-TocRun($paraIndex, $prefixText) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    [void]$rng.Find.Execute($prefixText, $true, $false, $false, $false, $false, $true, 1, $false, $prefixText, 2)
}

Merge-TocRun 5  "Assignment 01 "
Merge-TocRun 6  "Assignment 02  "
Merge-TocRun 7  "Assignment 03 "
Merge-TocRun 8  "Assignment 04   "
Merge-TocRun 13 "Assignment 01 "
Merge-TocRun 20 "Assignment 08 "
Merge-TocRun 22 "Assignment 01 "
Merge-TocRun 23 "Assignment 02 "

# ---------------------------------------------------------------------------
# 2) Small wording fixes.
# ---------------------------------------------------------------------------
$r1 = $d.Content
[void]$r1.Find.Execute("Write a code to print all the prime numbers that are present in it, using lambda expression.", $true, $false, $false, $false, $false, $true, 1, $false, "Write code to print all the prime numbers that are present in it, using lambda expression.", 2)

$r2 = $d.Content
[void]$r2.Find.Execute("Write acode to print all the Strings in reverse order, using lambda expression.", $true, $false, $false, $false, $false, $true, 1, $false, "Write code to print all the Strings in reverse order, using lambda expression.", 2)

# ---------------------------------------------------------------------------
# 3) Insert a new blank paragraph right after the "prime numbers" sentence
#    (matching the formatting of the blank paragraph that already follows
#    it, since InsertParagraphBefore clones the anchor paragraph's pPr/rPr).
# ---------------------------------------------------------------------------
$blankPara = $d.Paragraphs(34)
[void]$blankPara.Range.InsertParagraphBefore()

# ---------------------------------------------------------------------------
# 4) Bookmark "name" swaps: several headings carry two overlapping
#    bookmarks (a live TOC target + a stale duplicate); the edit swaps
#    which of the two names sits on the first vs. second bookmark.
# ---------------------------------------------------------------------------
function Swap-BookmarkNames($nameA, $nameB) {
    $bmA = $d.Bookmarks($nameA)
    $rangeA = $bmA.Range
    $bmA.Delete()

    $bmB = $d.Bookmarks($nameB)
    $rangeB = $bmB.Range
    $bmB.Delete()

    [void]$d.Bookmarks.Add($nameB, $rangeB)
    [void]$d.Bookmarks.Add($nameA, $rangeA)
}

Swap-BookmarkNames "_Toc4038"  "_Toc23397"
Swap-BookmarkNames "_Toc9697"  "_Toc12114"
Swap-BookmarkNames "_Toc1758"  "_Toc29283"
Swap-BookmarkNames "_Toc2345"  "_Toc31939"
Swap-BookmarkNames "_Toc9461"  "_Toc10089"
Swap-BookmarkNames "_Toc16655" "_Toc22551"

# ---------------------------------------------------------------------------
# 5) Move the "_GoBack" bookmark (last-edit marker) to sit right after the
#    "Example for perfect square numbers..." sentence.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

$anchor = $d.Content
[void]$anchor.Find.Execute("Example for perfect square numbers: 0, 1, 4, 9, 16, 25, 36, 49, 64 etc..", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)
[void]$d.Bookmarks.Add("_GoBack", $anchor)
